# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (F column) counters across the 展览 / 演出 / 本地生活 /
# 全部类型 sheets to the newly scraped values.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$values1 = @{
    2  = 21
    3  = 994
    4  = 238
    5  = 25
    6  = 1150
    7  = 928
    8  = 283
    10 = 81
    11 = 892
    12 = 318
    13 = 594
    14 = 525
    15 = 1374
    17 = 1268
    18 = 2931
    19 = 244
    20 = 1554
    21 = 1305
    22 = 753
    24 = 1306
    25 = 61
    28 = 3302
    29 = 644
    30 = 549
    31 = 1464
}
foreach ($row in $values1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $values1[$row]
}

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$values2 = @{
    4  = 14
    7  = 165
    10 = 4
    11 = 19
    14 = 8
}
foreach ($row in $values2.Keys) {
    $ws2.Cells.Item($row, 6).Value = $values2[$row]
}

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$values3 = @{
    2 = 774
}
foreach ($row in $values3.Keys) {
    $ws3.Cells.Item($row, 6).Value = $values3[$row]
}

# 全部类型 (All categories)
$ws4 = $wb.Worksheets.Item("全部类型")
$values4 = @{
    2  = 21
    3  = 774
    6  = 994
    7  = 238
    10 = 1150
    11 = 928
    12 = 283
    17 = 165
    21 = 4
    22 = 19
    23 = 892
    24 = 318
    25 = 594
    26 = 525
    27 = 1374
    29 = 1268
    30 = 2931
    31 = 244
    32 = 1554
    33 = 1305
    34 = 753
    36 = 1306
    42 = 3302
    43 = 644
    44 = 549
    45 = 1464
}
foreach ($row in $values4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $values4[$row]
}
